$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # 展览
$ws.Cells.Item(6, 6).Value2 = 1101
$ws.Cells.Item(8, 6).Value2 = 196
$ws.Cells.Item(9, 6).Value2 = 637
$ws.Cells.Item(11, 6).Value2 = 469
$ws.Cells.Item(12, 6).Value2 = 759
$ws.Cells.Item(13, 6).Value2 = 1458
$ws.Cells.Item(14, 6).Value2 = 1216
$ws.Cells.Item(15, 6).Value2 = 1445
$ws.Cells.Item(17, 6).Value2 = 1248
$ws.Cells.Item(18, 6).Value2 = 301
$ws.Cells.Item(19, 6).Value2 = 1590
$ws.Cells.Item(21, 6).Value2 = 1011
$ws.Cells.Item(22, 6).Value2 = 325
$ws.Cells.Item(25, 6).Value2 = 1388
$ws.Cells.Item(26, 6).Value2 = 85
$ws.Cells.Item(29, 6).Value2 = 1070
$ws.Cells.Item(30, 6).Value2 = 0
$ws.Cells.Item(31, 6).Value2 = 974
$ws.Cells.Item(32, 6).Value2 = 13
$ws.Cells.Item(34, 6).Value2 = 1319
$ws.Cells.Item(35, 6).Value2 = 1041
$ws.Cells.Item(37, 6).Value2 = 1068
$ws.Cells.Item(39, 6).Value2 = 49
$ws.Cells.Item(40, 6).Value2 = 40
$ws.Cells.Item(41, 6).Value2 = 834
$ws.Cells.Item(42, 6).Value2 = 1593
$ws.Cells.Item(44, 6).Value2 = 38
$ws.Cells.Item(45, 6).Value2 = 793
$ws.Cells.Item(48, 6).Value2 = 112
$ws = $wb.Worksheets.Item(2)  # 演出
$ws.Cells.Item(10, 6).Value2 = 172
$ws.Cells.Item(10, 7).Value2 = 280
$ws.Cells.Item(11, 6).Value2 = 1433
$ws.Cells.Item(12, 6).Value2 = 71
$ws.Cells.Item(14, 6).Value2 = 2532
$ws.Cells.Item(15, 6).Value2 = 1187
$ws.Cells.Item(16, 6).Value2 = 391
$ws.Cells.Item(17, 6).Value2 = 714
$ws.Cells.Item(18, 6).Value2 = 220
$ws.Cells.Item(20, 6).Value2 = 67
$ws.Cells.Item(23, 6).Value2 = 435
$ws.Cells.Item(26, 6).Value2 = 283
$ws.Cells.Item(27, 6).Value2 = 68326
$ws.Cells.Item(31, 6).Value2 = 181
$ws.Cells.Item(34, 6).Value2 = 136
$ws.Cells.Item(39, 6).Value2 = 166
$ws.Cells.Item(43, 6).Value2 = 39
$ws.Cells.Item(44, 6).Value2 = 39
$ws = $wb.Worksheets.Item(3)  # 本地生活
$ws.Cells.Item(4, 6).Value2 = 245
$ws.Cells.Item(5, 6).Value2 = 2771
$ws.Cells.Item(6, 6).Value2 = 4525
$ws.Cells.Item(10, 6).Value2 = 647
$ws.Cells.Item(11, 6).Value2 = 426
$ws.Cells.Item(12, 6).Value2 = 214
$ws.Cells.Item(13, 6).Value2 = 803
$ws.Cells.Item(14, 6).Value2 = 200
$ws.Cells.Item(15, 6).Value2 = 452
$ws = $wb.Worksheets.Item(4)  # 全部类型
$ws.Cells.Item(3, 6).Value2 = 245
$ws.Cells.Item(4, 6).Value2 = 2771
$ws.Cells.Item(5, 6).Value2 = 4525
$ws.Cells.Item(6, 6).Value2 = 647
$ws.Cells.Item(8, 6).Value2 = 214
$ws.Cells.Item(9, 6).Value2 = 214
$ws.Cells.Item(10, 6).Value2 = 803
$ws.Cells.Item(11, 6).Value2 = 803
$ws.Cells.Item(12, 6).Value2 = 200
$ws.Cells.Item(13, 6).Value2 = 1101
$ws.Cells.Item(15, 6).Value2 = 196
$ws.Cells.Item(16, 6).Value2 = 1433
$ws.Cells.Item(17, 6).Value2 = 469
$ws.Cells.Item(18, 6).Value2 = 759
$ws.Cells.Item(19, 6).Value2 = 2532
$ws.Cells.Item(20, 6).Value2 = 1187
$ws.Cells.Item(21, 6).Value2 = 1458
$ws.Cells.Item(22, 6).Value2 = 1216
$ws.Cells.Item(23, 6).Value2 = 1445
$ws.Cells.Item(24, 6).Value2 = 1248
$ws.Cells.Item(25, 6).Value2 = 220
$ws.Cells.Item(26, 6).Value2 = 67
$ws.Cells.Item(27, 6).Value2 = 1590
$ws.Cells.Item(29, 6).Value2 = 1011
$ws.Cells.Item(30, 6).Value2 = 325
$ws.Cells.Item(31, 6).Value2 = 452
$ws.Cells.Item(32, 6).Value2 = 452
$ws.Cells.Item(33, 6).Value2 = 435
$ws.Cells.Item(34, 6).Value2 = 1388
$ws.Cells.Item(37, 6).Value2 = 1070
$ws.Cells.Item(38, 6).Value2 = 283
$ws.Cells.Item(39, 6).Value2 = 974
$ws.Cells.Item(40, 6).Value2 = 13
$ws.Cells.Item(41, 6).Value2 = 1041
$ws.Cells.Item(43, 6).Value2 = 1068
$ws.Cells.Item(45, 6).Value2 = 835
$ws.Cells.Item(47, 6).Value2 = 1593
$ws.Cells.Item(49, 6).Value2 = 793
$ws.Cells.Item(50, 6).Value2 = 39
